$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.451.06'
$ws.Range("E2").Value = '  -1.10%  '

$ws.Range("D3").Value = '2.686.84'
$ws.Range("E3").Value = '  -2.79%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -2.77%  '

$ws.Range("E9").Value = '  -3.87%  '

$ws.Range("E10").Value = '  -0.96%  '

$ws.Range("E11").Value = '  -4.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.23%  '

$ws.Range("D13").Value = '3.159.34'
$ws.Range("E13").Value = '  -2.86%  '

$ws.Range("E14").Value = '  -1.97%  '

$ws.Range("D15").Value = '63.296.70'
$ws.Range("E15").Value = '  -0.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000146'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.01%  '

$ws.Range("D17").Value = '2.685.50'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.60%  '

$ws.Range("E22").Value = '  -0.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.509'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.27%  '

$ws.Range("D28").Value = '0.0₃0861'
$ws.Range("E28").Value = '  -5.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '165.87'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.83'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.70%  '

$ws.Range("E37").Value = '  -1.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '343.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.19%  '

$ws.Range("E39").Value = '  -6.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.16'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.13'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.620'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0564'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.40%  '

$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.84%  '

$ws.Range("E50").Value = '  -3.66%  '

$ws.Range("E51").Value = '  -4.55%  '
